# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (commit: "Updated cryptos list on Tue Jun 25 06:35:10 UTC
# 2024 with GitHub Actions"). Each entry maps an A1-style cell reference
# to its new text value. Rows 20 and 21 (Uniswap / Polkadot) also swap
# places in the ranking, so their Coin name, Link and Price/Volume cells
# are all rewritten.
#
# The "Price" column stores numbers as plain text (e.g. "571.82"), so for
# values that would otherwise be auto-recognized as a number by Excel we
# force the cell to Text format first to keep it a text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "61.065.09"; Numeric = $false },
    @{ Cell = "E2"; Value = "  -2.00%  "; Numeric = $false },
    @{ Cell = "D3"; Value = "3.372.21"; Numeric = $false },
    @{ Cell = "E3"; Value = "  -0.17%  "; Numeric = $false },
    @{ Cell = "E4"; Value = "  -0.02%  "; Numeric = $false },
    @{ Cell = "D5"; Value = "571.82"; Numeric = $true },
    @{ Cell = "D6"; Value = "135.78"; Numeric = $true },
    @{ Cell = "E6"; Value = "  +8.87%  "; Numeric = $false },
    @{ Cell = "D8"; Value = "3.371.50"; Numeric = $false },
    @{ Cell = "E8"; Value = "  -0.19%  "; Numeric = $false },
    @{ Cell = "E9"; Value = "  +0.73%  "; Numeric = $false },
    @{ Cell = "D10"; Value = "7.59"; Numeric = $true },
    @{ Cell = "E10"; Value = "  +4.78%  "; Numeric = $false },
    @{ Cell = "E11"; Value = "  +2.93%  "; Numeric = $false },
    @{ Cell = "E12"; Value = "  +4.07%  "; Numeric = $false },
    @{ Cell = "D13"; Value = "3.944.52"; Numeric = $false },
    @{ Cell = "E13"; Value = "  -0.21%  "; Numeric = $false },
    @{ Cell = "E14"; Value = "  +1.92%  "; Numeric = $false },
    @{ Cell = "E15"; Value = "  +1.57%  "; Numeric = $false },
    @{ Cell = "D16"; Value = "3.376.36"; Numeric = $false },
    @{ Cell = "E16"; Value = "  -0.13%  "; Numeric = $false },
    @{ Cell = "D17"; Value = "25.17"; Numeric = $true },
    @{ Cell = "E17"; Value = "  +3.31%  "; Numeric = $false },
    @{ Cell = "D18"; Value = "61.215.47"; Numeric = $false },
    @{ Cell = "E18"; Value = "  -1.87%  "; Numeric = $false },
    @{ Cell = "D19"; Value = "13.99"; Numeric = $true },
    @{ Cell = "E19"; Value = "  +7.37%  "; Numeric = $false },
    @{ Cell = "B20"; Value = "Polkadot"; Numeric = $false },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; Numeric = $false },
    @{ Cell = "D20"; Value = "5.81"; Numeric = $true },
    @{ Cell = "E20"; Value = "  +3.41%  "; Numeric = $false },
    @{ Cell = "B21"; Value = "Uniswap"; Numeric = $false },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; Numeric = $false },
    @{ Cell = "D21"; Value = "9.42"; Numeric = $true },
    @{ Cell = "E21"; Value = "  +2.46%  "; Numeric = $false },
    @{ Cell = "D22"; Value = "374.04"; Numeric = $true },
    @{ Cell = "E22"; Value = "  +1.01%  "; Numeric = $false },
    @{ Cell = "E23"; Value = "  +2.62%  "; Numeric = $false },
    @{ Cell = "D24"; Value = "3.505.65"; Numeric = $false },
    @{ Cell = "E24"; Value = "  -0.19%  "; Numeric = $false },
    @{ Cell = "E25"; Value = "  +0.11%  "; Numeric = $false },
    @{ Cell = "D26"; Value = "70.45"; Numeric = $true },
    @{ Cell = "E26"; Value = "  -0.62%  "; Numeric = $false },
    @{ Cell = "D27"; Value = "0.0000116"; Numeric = $true },
    @{ Cell = "E27"; Value = "  +10.95%  "; Numeric = $false },
    @{ Cell = "D28"; Value = "1.68"; Numeric = $true },
    @{ Cell = "E28"; Value = "  +22.41%  "; Numeric = $false },
    @{ Cell = "D29"; Value = "7.70"; Numeric = $true },
    @{ Cell = "E29"; Value = "  +12.44%  "; Numeric = $false },
    @{ Cell = "E30"; Value = "  -0.01%  "; Numeric = $false },
    @{ Cell = "D31"; Value = "8.12"; Numeric = $true },
    @{ Cell = "E31"; Value = "  +5.03%  "; Numeric = $false },
    @{ Cell = "E32"; Value = "  +2.36%  "; Numeric = $false },
    @{ Cell = "E33"; Value = "  +4.77%  "; Numeric = $false },
    @{ Cell = "E34"; Value = "  -0.04%  "; Numeric = $false },
    @{ Cell = "D35"; Value = "3.401.75"; Numeric = $false },
    @{ Cell = "E35"; Value = "  -0.15%  "; Numeric = $false },
    @{ Cell = "D36"; Value = "23.36"; Numeric = $true },
    @{ Cell = "E36"; Value = "  +3.22%  "; Numeric = $false },
    @{ Cell = "E37"; Value = "  +8.03%  "; Numeric = $false },
    @{ Cell = "E38"; Value = "  +5.27%  "; Numeric = $false },
    @{ Cell = "E39"; Value = "  +5.91%  "; Numeric = $false },
    @{ Cell = "D40"; Value = "163.79"; Numeric = $true },
    @{ Cell = "E40"; Value = "  -0.77%  "; Numeric = $false },
    @{ Cell = "D41"; Value = "0.0788"; Numeric = $true },
    @{ Cell = "E41"; Value = "  +5.36%  "; Numeric = $false },
    @{ Cell = "E42"; Value = "  -0.04%  "; Numeric = $false },
    @{ Cell = "E43"; Value = "  +4.40%  "; Numeric = $false },
    @{ Cell = "D44"; Value = "1.21"; Numeric = $true },
    @{ Cell = "E44"; Value = "  +13.11%  "; Numeric = $false },
    @{ Cell = "E45"; Value = "  -0.48%  "; Numeric = $false },
    @{ Cell = "D46"; Value = "41.29"; Numeric = $true },
    @{ Cell = "E47"; Value = "  +4.86%  "; Numeric = $false },
    @{ Cell = "D48"; Value = "23.19"; Numeric = $true },
    @{ Cell = "E48"; Value = "  +4.47%  "; Numeric = $false },
    @{ Cell = "D49"; Value = "6.98"; Numeric = $true },
    @{ Cell = "E49"; Value = "  +6.16%  "; Numeric = $false },
    @{ Cell = "D50"; Value = "22.82"; Numeric = $true },
    @{ Cell = "E50"; Value = "  +13.24%  "; Numeric = $false },
    @{ Cell = "D51"; Value = "0.893"; Numeric = $true },
    @{ Cell = "E51"; Value = "  +6.72%  "; Numeric = $false })

foreach ($u in $updates) {
    if ($u.Numeric) {
        $ws.Range($u.Cell).NumberFormat = "@"
    }
    $ws.Range($u.Cell).Value = $u.Value
}
